# update user & admin
# Replicates row 2 (the single sample "user" record) down through rows 3-9
# (adding 7 more "user/admin" records), wires up their e-mail hyperlinks,
# and updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate the data row (row 2) into rows 3-9, keeping styles ---
$ws.Range("A2:L2").Copy()
$ws.Range("A3:L3").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A4:L4").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A5:L5").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A6:L6").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A7:L7").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A8:L8").PasteSpecial()
$ws.Range("A2:L2").Copy()
$ws.Range("A9:L9").PasteSpecial()

# --- Fix up the incrementing "id" column ---
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8

# --- Re-create the e-mail hyperlinks for the new rows (order matches the
#     original authoring order so relationship ids line up: J3..J7, then J9, J8) ---
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J5"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J6"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J7"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J9"), "mailto:fadel.mm01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("J8"), "mailto:fadel.mm01@gmail.com")

# Adding a hyperlink re-styles the cell with the generic "Hyperlink" cell
# style; restore the original style (same as J2) so formatting matches.
$ws.Range("J3:J9").Style = $ws.Range("J2").Style

# --- Update the active selection shown when the workbook is opened ---
$ws.Range("D11:E11").Select()
